$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3
$ws.Range("L2").Value = 1.3
$ws.Range("M2").Value = 3.4
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 1.8
$ws.Range("T2").Value = 8.5
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 9.5
$ws.Range("W2").Value = 23
$ws.Range("X2").Value = 19
$ws.Range("Y2").Value = 29
$ws.Range("AE2").Value = 9.5
$ws.Range("AF2").Value = 15
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 29
$ws.Range("AI2").Value = 23

# Row 5
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 4.9
$ws.Range("I5").Value = 1.31
$ws.Range("T5").Value = 19.5
$ws.Range("V5").Value = 18.5
$ws.Range("Z5").Value = 17.5
$ws.Range("AB5").Value = 15
$ws.Range("AC5").Value = 50
$ws.Range("AE5").Value = 8
$ws.Range("AH5").Value = 7.6
$ws.Range("AI5").Value = 8.75
$ws.Range("AJ5").Value = 17.5

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 2.72
$ws.Range("T6").Value = 9.75
$ws.Range("U6").Value = 11.75
$ws.Range("V6").Value = 7.7
$ws.Range("W6").Value = 19.5
$ws.Range("X6").Value = 13
$ws.Range("AA6").Value = 6.5
$ws.Range("AE6").Value = 11.25
$ws.Range("AF6").Value = 15
$ws.Range("AG6").Value = 8.75
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 16.5

# Row 7
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 4.35
$ws.Range("V7").Value = 7
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 25
$ws.Range("AH7").Value = 60
$ws.Range("AI7").Value = 29

# Row 8
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("L8").Value = 1.25
$ws.Range("M8").Value = 3.75
$ws.Range("P8").Value = 1.36
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.95
$ws.Range("T8").Value = 7.5
$ws.Range("U8").Value = 8.5
$ws.Range("V8").Value = 8.5
$ws.Range("W8").Value = 15
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 26
$ws.Range("Z8").Value = 11
$ws.Range("AA8").Value = 6.5
$ws.Range("AB8").Value = 15
$ws.Range("AC8").Value = 51
$ws.Range("AD8").Value = 201
$ws.Range("AE8").Value = 13
$ws.Range("AG8").Value = 15
$ws.Range("AH8").Value = 51
$ws.Range("AI8").Value = 34
$ws.Range("AJ8").Value = 41

# Row 11
$ws.Range("L11").Value = 1.3
$ws.Range("M11").Value = 3.4

# Row 12
$ws.Range("G12").Value = 1.9
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 3.7
$ws.Range("N12").Value = 1.8
$ws.Range("O12").Value = 2
$ws.Range("R12").Value = 1.75
$ws.Range("S12").Value = 2
$ws.Range("AF12").Value = 21
$ws.Range("AG12").Value = 13

# Row 14
$ws.Range("G14").Value = 1.7
$ws.Range("H14").Value = 3.8
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 1.03
$ws.Range("K14").Value = 15
$ws.Range("L14").Value = 1.2
$ws.Range("M14").Value = 4.33
$ws.Range("N14").Value = 1.67
$ws.Range("O14").Value = 2.15
$ws.Range("P14").Value = 1.3
$ws.Range("Q14").Value = 3.4
$ws.Range("R14").Value = 1.67
$ws.Range("S14").Value = 2.1
$ws.Range("T14").Value = 8.5
$ws.Range("U14").Value = 9
$ws.Range("V14").Value = 8.5
$ws.Range("W14").Value = 15
$ws.Range("X14").Value = 13
$ws.Range("Y14").Value = 21
$ws.Range("Z14").Value = 13
$ws.Range("AA14").Value = 7.5
$ws.Range("AB14").Value = 13
$ws.Range("AC14").Value = 41
$ws.Range("AD14").Value = 151
$ws.Range("AE14").Value = 15
$ws.Range("AF14").Value = 23
$ws.Range("AG14").Value = 15
$ws.Range("AH14").Value = 51
$ws.Range("AI14").Value = 34
$ws.Range("AJ14").Value = 34
